$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = "MATTEO MARANER"
$ws.Range("B24").Value = "Stefano Tita | Clitoriders"
$ws.Range("C24").Value = "ANDREA MANFREDI | Pinguini Trentini"
$ws.Range("D24").Value = "Federico  Mortillaro | Clitoriders"
$ws.Range("E24").Value = "Federico  Manica | iMontagna"
$ws.Range("F24").Value = "Gentian Capa | FC. Stallions"
